# Add writing of a "DataType" column at the beginning of each row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
# Shift the existing headers one column to the right (A->B, B->C, C->D, D->E, E->F)
# and put the new "DataType" header in column A.
$ws.Range("F1").Value = "avgAgainst"
$ws.Range("E1").Value = "sumAgainst"
$ws.Range("D1").Value = "avgFor"
$ws.Range("C1").Value = "sumFor"
$ws.Range("B1").Value = "LearningRate"
$ws.Range("A1").Value = "DataType"

# Data type label repeating every three data rows.
$dataTypes = @("Before", "After", "Difference")

# New values for columns B (LearningRate), C (sumFor), D (avgFor),
# E (sumAgainst) -- column F stays empty for data rows.
$data = @(
    @(4330.81,  131.24,    4915.63,    148.96),
    @(4837.29,  146.58,    505.07,     15.31),
    @(506.48,   15.34,     -4410.56,   -133.65),
    @(4389.36,  133.01,    5074.63,    153.78),
    @(846.51,   25.65,     1050.26,    31.83),
    @(-3542.85, -107.36,   -4024.37,   -121.95),
    @(5003.86,  151.63,    4612.07,    139.76),
    @(2904.51,  88.02,     2895.74,    87.75),
    @(-2099.35, -63.61,    -1716.33,   -52.01),
    @(1389.59,  42.11,     2778.19,    84.19),
    @(2662.86,  80.69,     1736.37,    52.62),
    @(1273.27,  38.58,     -1041.82,   -31.57),
    @(759.74,   23.02,     839.87,     25.45),
    @(869.1900000000001, 26.34, 736.4, 22.32),
    @(109.45,   3.32,      -103.47,    -3.13)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $dataTypes[$i % 3]
    $ws.Range("B$row").Value = $data[$i][0]
    $ws.Range("C$row").Value = $data[$i][1]
    $ws.Range("D$row").Value = $data[$i][2]
    $ws.Range("E$row").Value = $data[$i][3]
}
